# Weekly update of the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Espinaca" data:
# each data row (2-24) is refreshed with the figures (date, volume, prices, origin)
# from the following week's record, rotating row 2's original values down to row 24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44798
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("O2").Value = "Provincia de Diguillín"
$ws.Range("P2").Value = 700
$ws.Range("D3").Value = 44211
$ws.Range("J3").Value = 28
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = 8214
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 821
$ws.Range("D4").Value = 44831
$ws.Range("J4").Value = 60
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 7250
$ws.Range("P4").Value = 725
$ws.Range("D5").Value = 44817
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("P5").Value = 700
$ws.Range("D6").Value = 44817
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 8000
$ws.Range("O6").Value = "Provincia de Diguillín"
$ws.Range("P6").Value = 800
$ws.Range("D7").Value = 44791
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 8500
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 8750
$ws.Range("P7").Value = 875
$ws.Range("D8").Value = 44819
$ws.Range("J8").Value = 100
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7500
$ws.Range("P8").Value = 750
$ws.Range("D9").Value = 44813
$ws.Range("J9").Value = 120
$ws.Range("D10").Value = 44841
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 6500
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 6750
$ws.Range("O10").Value = "Provincia de Diguillín"
$ws.Range("P10").Value = 675
$ws.Range("D11").Value = 44846
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 6500
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 6750
$ws.Range("P11").Value = 675
$ws.Range("D12").Value = 44782
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 8000
$ws.Range("M12").Value = 8500
$ws.Range("P12").Value = 850
$ws.Range("D13").Value = 44838
$ws.Range("K13").Value = 6500
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6750
$ws.Range("P13").Value = 675
$ws.Range("D14").Value = 44810
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("P14").Value = 750
$ws.Range("D15").Value = 44804
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7500
$ws.Range("M15").Value = 7250
$ws.Range("P15").Value = 725
$ws.Range("D16").Value = 44784
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 8000
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 8500
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 850
$ws.Range("D17").Value = 44775
$ws.Range("I17").Value = "Primera"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("D18").Value = 44806
$ws.Range("J18").Value = 120
$ws.Range("L18").Value = 7500
$ws.Range("M18").Value = 7250
$ws.Range("P18").Value = 725
$ws.Range("D19").Value = 44790
$ws.Range("K19").Value = 8500
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8750
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 875
$ws.Range("D20").Value = 44980
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7750
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 775
$ws.Range("D21").Value = 44812
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7500
$ws.Range("P21").Value = 750
$ws.Range("D22").Value = 44847
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 6500
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 6750
$ws.Range("P22").Value = 675
$ws.Range("D23").Value = 44799
$ws.Range("J23").Value = 60
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 7000
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 700
$ws.Range("D24").Value = 44203
$ws.Range("J24").Value = 27
$ws.Range("K24").Value = 7000
$ws.Range("M24").Value = 7556
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 756
